$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 13.109
$ws.Range("D3").Value = -7.726999999999999
$ws.Range("A4").Value = -22.09
$ws.Range("C4").Value = -12.927
$ws.Range("D4").Value = -7.727000000000001
$ws.Range("C5").Value = -12.78
$ws.Range("E5").Value = 12.886
$ws.Range("A6").Value = -21.345
$ws.Range("C6").Value = -12.712
$ws.Range("A7").Value = -21.118
$ws.Range("A8").Value = -21.584
$ws.Range("C8").Value = -12.653
$ws.Range("D9").Value = -8.048
$ws.Range("D11").Value = -7.8
$ws.Range("D14").Value = -8.028
$ws.Range("A16").Value = -21.04
$ws.Range("C16").Value = -12.862
$ws.Range("D18").Value = -7.637
$ws.Range("A20").Value = -22.145
$ws.Range("E20").Value = 13.229
$ws.Range("A21").Value = -20.921
$ws.Range("C22").Value = -12.78
$ws.Range("D25").Value = -7.913000000000001
